$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.412.30'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.191.36'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.27'
$ws.Range("E5").Value = '  +5.02%  '

$ws.Range("E6").Value = '  +0.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.72'
$ws.Range("E7").Value = '  -2.04%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +5.31%  '

$ws.Range("E10").Value = '  +1.95%  '

$ws.Range("E11").Value = '  +2.25%  '

$ws.Range("E12").Value = '  -1.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.17'
$ws.Range("E13").Value = '  +7.24%  '

$ws.Range("E14").Value = '  -0.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.517.12'
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("E16").Value = '  +4.81%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.55'
$ws.Range("E17").Value = '  -1.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.225.01'
$ws.Range("E18").Value = '  +1.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.303.70'
$ws.Range("E19").Value = '  +0.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0959'
$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("E21").Value = '  +1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.10'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.20'
$ws.Range("E23").Value = '  +0.85%  '

$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.96'
$ws.Range("E25").Value = '  +10.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.82'
$ws.Range("E26").Value = '  +21.27%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.53'
$ws.Range("E28").Value = '  +4.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.72'
$ws.Range("E29").Value = '  -5.11%  '

$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.81'
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("E32").Value = '  +1.49%  '

$ws.Range("E33").Value = '  -0.61%  '

$ws.Range("E34").Value = '  +7.36%  '

$ws.Range("E35").Value = '  -0.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.50'
$ws.Range("E36").Value = '  +6.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.83'
$ws.Range("E37").Value = '  +12.36%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.16'
$ws.Range("E38").Value = '  +7.73%  '

$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.64'
$ws.Range("E39").Value = '  +0.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0299'
$ws.Range("E40").Value = '  +8.67%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("E41").Value = '  -2.79%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.49'
$ws.Range("E42").Value = '  +18.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.70'
$ws.Range("E43").Value = '  -2.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.55'
$ws.Range("E44").Value = '  +0.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.13'
$ws.Range("E45").Value = '  +4.48%  '

$ws.Range("E46").Value = '  +3.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.64'
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.102'
$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  +7.77%  '

$ws.Range("E51").Value = '  +0.71%  '
